# Update the repo URL placeholder in the "Clone the Repository" code block.
#
# Before: git clone <your-repo-url>
# After:  git clone https://github.com/ml-multimedia-hit-2025/deepcompress.git
#
# The original text is spread across four runs with different character
# styles ( NormalTok " clone ", OperatorTok "<", NormalTok "your-repo-url",
# OperatorTok ">" ). A single Find/Replace across that span collapses it
# into the first run's style (NormalTok) and removes the rest, which is
# exactly what the target diff shows.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    " clone <your-repo-url>",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    " clone https://github.com/ml-multimedia-hit-2025/deepcompress.git",
    2
)
